$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 422, shifting existing rows 422-508 down to 423-509
$ws.Rows(422).Insert()

# Populate the newly inserted row 422 with the new data record
$ws.Range("A422").Value = 3
$ws.Range("B422").Value = "Femacal de La Calera"
$ws.Range("C422").Value = "Coquimbo"
$ws.Range("D422").Value = 45015
$ws.Range("E422").Value = 5
$ws.Range("F422").Value = 100112012
$ws.Range("G422").Value = "Espinaca"
$ws.Range("H422").Value = "Sin especificar"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 150
$ws.Range("K422").Value = 5500
$ws.Range("L422").Value = 6000
$ws.Range("M422").Value = 5767
$ws.Range("N422").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O422").Value = "Provincia de Quillota"
$ws.Range("P422").Value = 1922
$ws.Range("Q422").Value = 3
$ws.Range("R422").Value = "Hortaliza"
